$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 48.77778
$ws.Range("I4").Value = 48.77778
$ws.Range("K4").Value = 48.77778
$ws.Range("M4").Value = 65.22221999999999

$ws.Range("H129").Value = 2308.6667
$ws.Range("I129").Value = 1730.3334
$ws.Range("J129").Value = 2887
$ws.Range("K129").Value = 5191.0002
$ws.Range("L129").Value = 8661
$ws.Range("M129").Value = -191.0002000000004
$ws.Range("N129").Value = -18661

$ws.Range("H132").Value = 1160.3214
$ws.Range("I132").Value = 1181.3091
$ws.Range("K132").Value = 3543.9273
$ws.Range("M132").Value = -1013.9273

$ws.Range("H137").Value = 5142.9375
$ws.Range("I137").Value = 3686.625
$ws.Range("K137").Value = 11059.875
$ws.Range("M137").Value = -8509.875

$ws.Range("H138").Value = 6812.0635
$ws.Range("I138").Value = 5080
$ws.Range("J138").Value = 6961.3794
$ws.Range("K138").Value = 15240
$ws.Range("L138").Value = 20884.1382
$ws.Range("M138").Value = -10100
$ws.Range("N138").Value = -31164.1382

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19310.033
$ws.Range("I32").Value = 11594.637
$ws.Range("K32").Value = 11594.637
$ws.Range("M32").Value = -11307.637

$ws.Range("H61").Value = 113455
$ws.Range("I61").Value = 2297.4285
$ws.Range("K61").Value = 2297.4285
$ws.Range("M61").Value = -2085.4285

$ws.Range("H136").Value = 113455
$ws.Range("I136").Value = 2297.4285
$ws.Range("K136").Value = 6892.2855
$ws.Range("M136").Value = -4342.2855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7350
$ws.Range("I86").Value = 4800
$ws.Range("J86").Value = 9900
$ws.Range("K86").Value = 4800
$ws.Range("L86").Value = 9900
$ws.Range("M86").Value = -3677
$ws.Range("N86").Value = -12146

$ws.Range("H89").Value = 7350
$ws.Range("I89").Value = 4800
$ws.Range("J89").Value = 9900
$ws.Range("K89").Value = 24000
$ws.Range("L89").Value = 49500
$ws.Range("M89").Value = -18384
$ws.Range("N89").Value = -60732

$ws.Range("H134").Value = 2484.074
$ws.Range("I134").Value = 2116.8635
$ws.Range("K134").Value = 6350.5905
$ws.Range("M134").Value = -3815.5905

$ws.Range("H141").Value = 186165.67
$ws.Range("J141").Value = 186165.67
$ws.Range("L141").Value = 186165.67
$ws.Range("N141").Value = -196525.67

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2050.75
$ws.Range("I58").Value = 1944.5714
$ws.Range("J58").Value = 2199.4
$ws.Range("K58").Value = 1944.5714
$ws.Range("L58").Value = 2199.4
$ws.Range("M58").Value = -1741.5714
$ws.Range("N58").Value = -2605.4

$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws.Range("H131").Value = 59996
$ws.Range("I131").Value = 30000
$ws.Range("J131").Value = 69994.664
$ws.Range("K131").Value = 30000
$ws.Range("L131").Value = 69994.664
$ws.Range("M131").Value = -24960
$ws.Range("N131").Value = -80074.664

$ws.Range("H136").Value = 2050.75
$ws.Range("I136").Value = 1944.5714
$ws.Range("J136").Value = 2199.4
$ws.Range("K136").Value = 5833.7142
$ws.Range("L136").Value = 6598.200000000001
$ws.Range("M136").Value = -3283.7142
$ws.Range("N136").Value = -11698.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 913.3333
$ws.Range("I97").Value = 881.1539
$ws.Range("K97").Value = 881.1539
$ws.Range("M97").Value = -385.1539

$ws.Range("H109").Value = 87001
$ws.Range("J109").Value = 87001
$ws.Range("L109").Value = 87001
$ws.Range("N109").Value = -89081

$ws.Range("H110").Value = 98854.75
$ws.Range("J110").Value = 98854.75
$ws.Range("L110").Value = 98854.75
$ws.Range("N110").Value = -107034.75

$ws.Range("H119").Value = 66843.14
$ws.Range("J119").Value = 69704
$ws.Range("L119").Value = 69704
$ws.Range("N119").Value = -79380

$ws.Range("H132").Value = 6231.048
$ws.Range("I132").Value = 6211.8823
$ws.Range("J132").Value = 6312.5
$ws.Range("K132").Value = 18635.6469
$ws.Range("L132").Value = 18937.5
$ws.Range("M132").Value = -16105.6469
$ws.Range("N132").Value = -23997.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1170
$ws.Range("I22").Value = 340
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 340
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = -45
$ws.Range("N22").Value = -2590

$ws.Range("H27").Value = 1170
$ws.Range("I27").Value = 340
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 340
$ws.Range("L27").Value = 2000
$ws.Range("M27").Value = -233
$ws.Range("N27").Value = -2214

$ws.Range("H46").Value = 2730
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 2730
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 2730
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -3106

$ws.Range("H55").Value = 930.2973
$ws.Range("I55").Value = 511.25
$ws.Range("K55").Value = 511.25
$ws.Range("M55").Value = -338.25

$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws.Range("H133").Value = 74375.336
$ws.Range("J133").Value = 74375.336
$ws.Range("L133").Value = 74375.336
$ws.Range("N133").Value = -79435.336

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7500
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 7500
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H100").Value = 1258
$ws.Range("I100").Value = 1400
$ws.Range("K100").Value = 2800
$ws.Range("M100").Value = -2259

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

$ws.Range("H126").Value = 9221.223
$ws.Range("I126").Value = 1879.2174
$ws.Range("K126").Value = 5637.6522
$ws.Range("M126").Value = -3167.6522

$ws.Range("H132").Value = 1259.8889
$ws.Range("I132").Value = 1062.8572
$ws.Range("K132").Value = 3188.5716
$ws.Range("M132").Value = -658.5715999999998

$ws.Range("H136").Value = 4347.8184
$ws.Range("I136").Value = 3502.389
$ws.Range("J136").Value = 5362.3335
$ws.Range("K136").Value = 10507.167
$ws.Range("L136").Value = 16087.0005
$ws.Range("M136").Value = -7957.167000000001
$ws.Range("N136").Value = -21187.0005
